$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 406 (existing rows 406..467 shift down to 407..468)
$ws.Rows("406:406").Insert()

# Populate the newly inserted row 406 with the new data record
$ws.Range("A406").Value = 5
$ws.Range("B406").Value = "Macroferia Regional de Talca"
$ws.Range("C406").Value = "Maule"
$ws.Range("D406").Value = 45127
$ws.Range("E406").Value = 7
$ws.Range("F406").Value = "Fruta"
$ws.Range("G406").Value = 100101
$ws.Range("H406").Value = "Berries"
$ws.Range("I406").Value = 100101007
$ws.Range("J406").Value = "Kiwi"
$ws.Range("K406").Value = "Hayward"
$ws.Range("L406").Value = "Primera"
$ws.Range("M406").Value = 220
$ws.Range("N406").Value = 14000
$ws.Range("O406").Value = 14000
$ws.Range("P406").Value = 14000
$ws.Range("Q406").Value = "$/bandeja 18 kilos"
$ws.Range("R406").Value = "Provincia de Curicó"
$ws.Range("S406").Value = 778
$ws.Range("T406").Value = 18
